# Weekly update: a new price-report row for "Feria Lagunitas de Puerto
# Montt" / Chirimoya is inserted at row 10, pushing the existing rows
# 10-22 down to 11-23 (dimension grows from A1:T22 to A1:T23).
#
# The new row 10 repeats the fixed/reference columns (Mercado ID, Mercado,
# Región, Codreg, Tipo, Producto ID, Producto, Categoría ID, Categoría,
# Variedad, Unidad de comercialización, Origen, Kg/unidad) from the row
# immediately below it (the row that used to be row 10), and carries its
# own new Fecha/Calidad/Volumen/Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 10..22 down to 11..23, leaving a blank row 10 in place.
$ws.Rows.Item(10).Insert()

# Columns that stay identical to the row below (now row 11, formerly row 10).
$ws.Cells.Item(10, 1).Value  = $ws.Cells.Item(11, 1).Value2   # Mercado ID
$ws.Cells.Item(10, 2).Value  = $ws.Cells.Item(11, 2).Value2   # Mercado
$ws.Cells.Item(10, 3).Value  = $ws.Cells.Item(11, 3).Value2   # Región
$ws.Cells.Item(10, 5).Value  = $ws.Cells.Item(11, 5).Value2   # Codreg
$ws.Cells.Item(10, 6).Value  = $ws.Cells.Item(11, 6).Value2   # Tipo
$ws.Cells.Item(10, 7).Value  = $ws.Cells.Item(11, 7).Value2   # Producto ID
$ws.Cells.Item(10, 8).Value  = $ws.Cells.Item(11, 8).Value2   # Producto
$ws.Cells.Item(10, 9).Value  = $ws.Cells.Item(11, 9).Value2   # Categoría ID
$ws.Cells.Item(10, 10).Value = $ws.Cells.Item(11, 10).Value2  # Categoría
$ws.Cells.Item(10, 11).Value = $ws.Cells.Item(11, 11).Value2  # Variedad
$ws.Cells.Item(10, 17).Value = $ws.Cells.Item(11, 17).Value2  # Unidad de comercialización
$ws.Cells.Item(10, 18).Value = $ws.Cells.Item(11, 18).Value2  # Origen
$ws.Cells.Item(10, 20).Value = $ws.Cells.Item(11, 20).Value2  # Kg / unidad

# New data for the inserted row.
$ws.Cells.Item(10, 4).Value  = 44880     # Fecha
$ws.Cells.Item(10, 12).Value = "Primera" # Calidad
$ws.Cells.Item(10, 13).Value = 300       # Volumen
$ws.Cells.Item(10, 14).Value = 22000     # Precio mínimo
$ws.Cells.Item(10, 15).Value = 22500     # Precio máximo
$ws.Cells.Item(10, 16).Value = 22250     # Precio promedio ponderado
$ws.Cells.Item(10, 19).Value = 2781      # Precio $/Kg
